$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trial input variables (B3, D3, E3, F3) - formulas in G3:K3 auto-recalculate.
$ws.Range("B3").Value = 330
$ws.Range("D3").Value = 180
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 60

# Drop the stray "protected, no fill" cell style from the unfilled label/result
# cells (row 2 labels + row 3 results) -- they fall back to the plain default
# style, same as Excel does when it re-saves and consolidates the cellXfs
# table.
$ws.Range("A2:F2").Style = "Normal"
$ws.Range("G3:K3").Style = "Normal"

# Move the active selection to A3 (matches the saved sheet view state).
$ws.Range("A3").Select() | Out-Null
